# datos_de_monedas.xlsx — refresh the coin-market snapshot numbers.
# Columns E (24h %) and F (7d %) are stored as text even though they look
# like numbers, so a leading apostrophe is used to keep Excel from
# re-typing them as numeric literals. Column G (price in btc) is the same
# text-number convention but only some rows changed. Column H (coins
# comercializadas) is a genuine numeric column.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - BTC
$ws.Range("E2").Value = "'-1.61"
$ws.Range("F2").Value = "'4.28"
$ws.Range("H2").Value = 27432954554.97071

# Row 3 - ETH
$ws.Range("E3").Value = "'-1.53"
$ws.Range("F3").Value = "'1.21"
$ws.Range("G3").Value = "'0.048538"
$ws.Range("H3").Value = 10698834803.34144

# Row 4 - USDT
$ws.Range("E4").Value = "'-0.03"
$ws.Range("F4").Value = "'-0.03"
$ws.Range("H4").Value = 42317824494.95636

# Row 5 - BNB
$ws.Range("E5").Value = "'-1.09"
$ws.Range("F5").Value = "'2.74"
$ws.Range("G5").Value = "'0.009398"
$ws.Range("H5").Value = 518175263.1813599

# Row 6 - SOL
$ws.Range("E6").Value = "'-3.87"
$ws.Range("F6").Value = "'18.50"
$ws.Range("G6").Value = "'0.002371"
$ws.Range("H6").Value = 2579851111.428535

# Row 7 - USDC (only the coins-comercializadas number moved)
$ws.Range("H7").Value = 3499145249.343101

# Row 8 - STETH
$ws.Range("E8").Value = "'-1.55"
$ws.Range("F8").Value = "'1.16"
$ws.Range("G8").Value = "'0.048478"
$ws.Range("H8").Value = 86071853.94508265

# Row 9 - XRP
$ws.Range("E9").Value = "'-2.47"
$ws.Range("F9").Value = "'5.36"
$ws.Range("H9").Value = 1573477918.634741

# Row 10 - DOGE
$ws.Range("E10").Value = "'-4.15"
$ws.Range("F10").Value = "'15.70"
$ws.Range("H10").Value = 1115770111.721719

# Row 11 - TON
$ws.Range("E11").Value = "'-1.37"
$ws.Range("F11").Value = "'14.28"
$ws.Range("G11").Value = "'0.000094"
$ws.Range("H11").Value = 108421926.3558878

# Row 12 - ADA
$ws.Range("E12").Value = "'-1.37"
$ws.Range("F12").Value = "'1.90"
$ws.Range("H12").Value = 282407397.525009
